# Finanzen workbook update:
#  1. Add a new row to "Kredite" for an additional competitor credit offer ("Sofortkredit").
#  2. Add a new "Transportkosten" sheet listing shipping/transport cost options (competitors).

$wb = $excel.ActiveWorkbook

# --- 1. Kredite: add the "Sofortkredit" competitor row -----------------
$kredite = $wb.Worksheets.Item("Kredite")
$kredite.Range("A5").Value = "Sofortkredit"
$kredite.Range("B5").Value = 15000
$kredite.Range("C5").Value = 12.9
$kredite.Range("D5").Value = 12

# --- 2. New "Transportkosten" sheet -------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Transportkosten"

# Reuse the same bold/bordered/centered header style used on the other sheets
$srcHeader = $kredite.Range("A1:D1")
$dstHeader = $ws.Range("A1:D1")
$srcHeader.Copy()
$dstHeader.PasteSpecial(-4122)

# Header row
$ws.Range("A1").Value = "Transportart"
$ws.Range("B1").Value = "Kosten_pro_km"
$ws.Range("C1").Value = "Basis_Transportkosten"
$ws.Range("D1").Value = "Mindestkosten"

# Data rows
$ws.Range("A2").Value = "Standard Lieferung"
$ws.Range("B2").Value = 0.5
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 10

$ws.Range("A3").Value = "Express Lieferung"
$ws.Range("B3").Value = 0.8
$ws.Range("C3").Value = 8
$ws.Range("D3").Value = 15

$ws.Range("A4").Value = "Sperrgut"
$ws.Range("B4").Value = 1.2
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 25
